$d = $word.ActiveDocument
# Try using Find/Replace (pure text) only, skip InsertXML entirely, to establish baseline with ZERO drift expected
$d.Content.Find.Execute("zzzznonexistent", $true, $false, $false, $false, $false, $true, 1, $false, "zz", 2)
